$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 60856.58780776655
    3 = 57881.17288160624
    4 = 55399.93828265802
    5 = 53707.54692618801
    6 = 53680.8856395722
    7 = 56634.9832135402
    8 = 58446.84605949695
    9 = 63142.72549009701
    10 = 76716.99583502614
    11 = 83491.60492431262
    12 = 87563.66895977985
    13 = 90303.87260687441
    14 = 90165.42497353299
    15 = 93864.75891707004
    16 = 94482.71627596006
    17 = 92926.86975363505
    18 = 89804.61302349941
    19 = 83681.57571663187
    20 = 83714.51021982718
    21 = 79427.54689600211
    22 = 76549.06870090037
    23 = 74386.40041606507
    24 = 71540.08295675143
    25 = 67740.71402363165

}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $values[$row]
}
